$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 38), columns A (=1) through Y (=25)
$rowValues = @(36, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 430, 102, 1456, 1988, 0, 0, 0, 3, 2, 0, 66)

$newRow = 38
$lastRow = 37

for ($col = 1; $col -le $rowValues.Length; $col++) {
    $ws.Cells.Item($newRow, $col).Value = $rowValues[$col - 1]
}

# Copy the formatting of column A from the previous data row (bold, centered, bordered)
# so the new row's index cell matches the rest of the table.
$ws.Cells.Item($lastRow, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)

$excel.CutCopyMode = $false
